# Update the "Senior Computer Scientist - Web Standards" job-opening bullet
# on the single slide so it reads "Senior Computer Scientist - Web Platform
# Innovation and Standards", split across three runs the way the authored
# edit does:
#   "Senior Computer Scientist - "  (existing run, text shortened)
#   "Web Platform Innovation "      (new run)
#   "and Standards"                 (new run)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldParaText = "Senior Computer Scientist – Web Standards"
$oldPhrase   = "Web Standards"
$newPhrase   = "Web Platform Innovation and Standards"
$tailPhrase  = "and Standards"

# Find the shape on the slide whose text contains the bullet we need to edit.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "*$oldParaText*") {
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count

# Find the exact paragraph that holds the bullet (trim the trailing
# paragraph-mark character before comparing text).
$targetPara = $null
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    $candidateText = $candidate.Text.TrimEnd("`r")
    if ($candidateText -eq $oldParaText) {
        $targetPara = $candidate
    }
}

# Replace "Web Standards" with the new phrase; this keeps the leading
# "Senior Computer Scientist - " run untouched and creates a new run for
# the replacement text.
$currentText = $targetPara.Text
$phrasePos = $currentText.IndexOf($oldPhrase)
$replaceRange = $targetPara.Characters($phrasePos + 1, $oldPhrase.Length)
$replaceRange.Text = $newPhrase

# Split "and Standards" back out into its own run so the paragraph ends up
# as three runs, matching the authored edit.
$currentText = $targetPara.Text
$tailPos = $currentText.IndexOf($tailPhrase)
$tailRange = $targetPara.Characters($tailPos + 1, $tailPhrase.Length)
$tailRange.Text = $tailPhrase
